# removed annoying ! from calc field
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "U/I "
$ws.Range("C2").Value = "Fr(2,val1)+Fr(val3,val2) "
$ws.Range("E2").Value = "val1+val2 "
$ws.Range("G2").Value = "F*G"
$ws.Range("A3").Value = "U**2/R "
$ws.Range("D3").Value = '["Ein Quadrat hat eine Seiten länge von $val1 cm", ["val1*4", "Wie groß ist der Umfang"], ["val1**2", "Wie groß ist die Fläche"], ["val2*4","Wie groß wäre der Umfang, wenn die Seitenlänge $val2 cm wäre."]] '
$ws.Range("A4").Value = "R*I "
$ws.Range("A5").Value = "U/(R*1000) "
$ws.Range("A6").Value = "U/(R*1000) "

$ws.Range("D3").Select()
